# Add a new date column (CH) for 2024/12/03 to the "合成確率" sheet.
# CG (column 85) was the last existing date column (2024/12/02); this adds
# CH (column 86) with header text "2024/12/03" plus 52 numeric data points
# (rows 2-53), reusing the same three cell styles already used throughout
# the sheet (1 = plain, 2 = yellow "low" highlight, 3 = blue "mid" highlight).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Style-index -> a pre-existing cell carrying that exact style, used purely
# as a format donor for PasteSpecial (xlPasteFormats = -4122) so the new
# column reuses the workbook's existing style entries instead of minting
# new ones.
$styleTemplate = @{ 1 = "A2"; 2 = "D2"; 3 = "N2" }

# Give column CH the same rendered width as the rest of the table.
# (Raw OOXML <col width="12"/> corresponds to ColumnWidth 11.17 in the
# Excel object model's measurement units.)
$ws.Range("CH1").ColumnWidth = 11.17

# ---- Header cell CH1 -----------------------------------------------------
# Force text storage first (otherwise "2024/12/03" gets auto-recognised as
# a date serial), then restyle to match the rest of the header row.
$ws.Range("CH1").NumberFormat = "@"
$ws.Range("CH1").Value = "2024/12/03"
$ws.Range("CG1").Copy()
$ws.Range("CH1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Data cells CH2:CH53 --------------------------------------------------
$styles = @(2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,3,3,1,3,1,1,3,1,1,1,1,3,1,1,1,3,1,1,1,1,1,2,3,1,1,1,1,1,1,1,2,3,1)
$values = @(120.4,151.7,148.3,278.2,211,165.7,145.4,154.3,226.2,174.7,141.6,173.3,195.2,140.8,180.9,153.8,142.8,188.9,196.2,128.5,139,189.2,133.9,214.6,225.4,139,353.6,163.4,206.6,204.1,125.9,148.9,187.2,266.3,139.9,148,155.4,143.1,151.6,161.4,123.6,127.3,233.8,178.5,142.6,162.2,196.1,190.7,150.5,118.6,134.3,156.9)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $style = $styles[$i]
    $template = $styleTemplate[$style]

    $ws.Range($template).Copy()
    $ws.Range("CH$row").PasteSpecial(-4122)
    $ws.Range("CH$row").Value = $values[$i]
}
$excel.CutCopyMode = 0
